$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header label (Input inside body -> Method/Input inside body) ---
$ws.Range("B1").Value = "Method/Input inside body"

# --- "/" endpoint: add GET method ---
$ws.Range("B3").Value = "get"

# --- "/api" endpoint: add GET method + adjust return payload shape ---
$ws.Range("B5").Value = "get"
$ws.Range("C5").Value = "{data:{message:”/api accessible”}}"

# --- "/api/login" endpoint: add POST body spec ---
$ws.Range("B7").Value = "post/{username,password}"

# --- "/api/rooms" endpoint: add GET method ---
$ws.Range("B9").Value = "get"

# --- "/api/rooms/:id" endpoint: add GET method ---
$ws.Range("B10").Value = "get"

# --- New row: "/api/rooms/:id" DELETE method (inserted logically after the GET row) ---
$ws.Range("A11").Value = "/api/rooms/:id"
$ws.Range("B11").Value = "delete"
$ws.Range("C11").Value = "{data:{message:”success”}}"
$ws.Range("D11").Value = "note the special delete method"

# --- "/api/rooms/create" endpoint shifts down one row, gets POST body spec ---
$ws.Range("A12").Value = "/api/rooms/create"
$ws.Range("B12").Value = "post/{number,baseRent}"
$ws.Range("C12").Value = "{data:{room:{}}}"
$ws.Range("D12").Value = "The newly created Room"

# --- Row 13 becomes the blank separator row ---
$ws.Range("A13").Value = ""
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""

# --- "/api/tenants" endpoint shifts down one row ---
$ws.Range("A14").Value = "/api/tenants"
$ws.Range("C14").Value = "{data:{tenants:[]}}"

# --- "/api/tenants/:id" endpoint shifts down one row ---
$ws.Range("A15").Value = "/api/tenants/:id"
$ws.Range("B15").Value = ""
$ws.Range("D15").Value = ""

# --- New row: "/api/tenants/:id" DELETE method ---
$ws.Range("A16").Value = "/api/tenants/:id"
$ws.Range("B16").Value = "delete"
$ws.Range("C16").Value = "{data:{message:”success”}}"

# --- "/api/tenants/create" endpoint shifts down two rows total ---
$ws.Range("A17").Value = "/api/tenants/create"
$ws.Range("B17").Value = "{name,phoneNumber,aadharCard,room}"
$ws.Range("C17").Value = "{data:{tenant:{}}}"
$ws.Range("D17").Value = "The newly created Tenant"

# --- Row 18 becomes the blank separator row ---
$ws.Range("A18").Value = ""
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""

# --- "/api/transactions/?room" endpoint shifts down two rows, gets GET method ---
$ws.Range("A19").Value = "/api/transactions/?room"
$ws.Range("B19").Value = "get"
$ws.Range("C19").Value = "{data:{transactions:[]}}"
$ws.Range("D19").Value = ""

# --- "/api/transactions/:id" endpoint shifts down two rows, gets GET method ---
$ws.Range("A20").Value = "/api/transactions/:id"
$ws.Range("B20").Value = "get"
$ws.Range("C20").Value = "{data:{transaction:{}}}"

# --- "/api/tenants/create/?roomNumber" endpoint shifts down two rows, gets POST body spec ---
$ws.Range("A21").Value = "/api/tenants/create/?roomNumber"
$ws.Range("B21").Value = "post/{room,balance,transfer,remarks}"
$ws.Range("C21").Value = "{data:{transaction:{}}}"
$ws.Range("D21").Value = "The newly created Transaction"

# --- Extend the used range by two blank, identically formatted rows (35 & 36) ---
$ws.Range("A34:AI34").Copy()
$ws.Range("A35:AI36").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Move the active selection to D6, matching the saved view state ---
$ws.Range("D6").Select()
